$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 85,7

$data[0,0] = "Politik"
$data[0,1] = "https://www.cnnindonesia.com/nasional/politik"
$data[0,2] = "purbaya"
$data[0,3] = 1
$data[0,4] = 45931.37121663953
$data[0,5] = $null
$data[0,6] = "2025-10-01 08:54:33.117656+07:00"

$data[1,0] = "Hukum & Kriminal"
$data[1,1] = "https://www.cnnindonesia.com/nasional/hukum-kriminal"
$data[1,2] = "purbaya"
$data[1,3] = 1
$data[1,4] = 45931.37121804828
$data[1,5] = $null
$data[1,6] = "2025-10-01 08:54:33.239371+07:00"

$data[2,0] = "Peristiwa"
$data[2,1] = "https://www.cnnindonesia.com/nasional/peristiwa"
$data[2,2] = "purbaya"
$data[2,3] = 1
$data[2,4] = 45931.37121943914
$data[2,5] = $null
$data[2,6] = "2025-10-01 08:54:33.359542+07:00"

$data[3,0] = "Pemilu"
$data[3,1] = "https://www.cnnindonesia.com/nasional/pemilu"
$data[3,2] = "purbaya"
$data[3,3] = 1
$data[3,4] = 45931.37122100837
$data[3,5] = $null
$data[3,6] = "2025-10-01 08:54:33.495123+07:00"

$data[4,0] = "Info Politik"
$data[4,1] = "https://www.cnnindonesia.com/nasional/info-politik"
$data[4,2] = "purbaya"
$data[4,3] = 1
$data[4,4] = 45931.37122244346
$data[4,5] = $null
$data[4,6] = "2025-10-01 08:54:33.619115+07:00"

$data[5,0] = "Asia Pasifik"
$data[5,1] = "https://www.cnnindonesia.com/internasional/asia-pasifik"
$data[5,2] = "purbaya"
$data[5,3] = 1
$data[5,4] = 45931.37122411404
$data[5,5] = $null
$data[5,6] = "2025-10-01 08:54:33.763453+07:00"

$data[6,0] = "Timur Tengah"
$data[6,1] = "https://www.cnnindonesia.com/internasional/timur-tengah"
$data[6,2] = "purbaya"
$data[6,3] = 1
$data[6,4] = 45931.37122486073
$data[6,5] = $null
$data[6,6] = "2025-10-01 08:54:33.827967+07:00"

$data[7,0] = "Eropa Amerika"
$data[7,1] = "https://www.cnnindonesia.com/internasional/eropa-amerika"
$data[7,2] = "purbaya"
$data[7,3] = 1
$data[7,4] = 45931.37122617446
$data[7,5] = $null
$data[7,6] = "2025-10-01 08:54:33.941473+07:00"

$data[8,0] = "Keuangan"
$data[8,1] = "https://www.cnnindonesia.com/ekonomi/keuangan"
$data[8,2] = "purbaya"
$data[8,3] = 1
$data[8,4] = 45931.37122703999
$data[8,5] = $null
$data[8,6] = "2025-10-01 08:54:34.016255+07:00"

$data[9,0] = "Energi"
$data[9,1] = "https://www.cnnindonesia.com/ekonomi/energi"
$data[9,2] = "purbaya"
$data[9,3] = 1
$data[9,4] = 45931.37122847586
$data[9,5] = $null
$data[9,6] = "2025-10-01 08:54:34.140314+07:00"

$data[10,0] = "Bisnis"
$data[10,1] = "https://www.cnnindonesia.com/ekonomi/bisnis"
$data[10,2] = "purbaya"
$data[10,3] = 1
$data[10,4] = 45931.37123021537
$data[10,5] = $null
$data[10,6] = "2025-10-01 08:54:34.290608+07:00"

$data[11,0] = "Corporate Action"
$data[11,1] = "https://www.cnnindonesia.com/ekonomi/corporate-action"
$data[11,2] = "purbaya"
$data[11,3] = 1
$data[11,4] = 45931.37123152721
$data[11,5] = $null
$data[11,6] = "2025-10-01 08:54:34.403951+07:00"

$data[12,0] = "Sepakbola"
$data[12,1] = "https://www.cnnindonesia.com/olahraga/sepakbola"
$data[12,2] = "purbaya"
$data[12,3] = 1
$data[12,4] = 45931.37123255712
$data[12,5] = $null
$data[12,6] = "2025-10-01 08:54:34.492935+07:00"

$data[13,0] = "Moto GP"
$data[13,1] = "https://www.cnnindonesia.com/olahraga/moto-gp"
$data[13,2] = "purbaya"
$data[13,3] = 1
$data[13,4] = 45931.37123433159
$data[13,5] = $null
$data[13,6] = "2025-10-01 08:54:34.646249+07:00"

$data[14,0] = "Teknologi Informasi"
$data[14,1] = "https://www.cnnindonesia.com/teknologi/teknologi-informasi"
$data[14,2] = "purbaya"
$data[14,3] = 1
$data[14,4] = 45931.37123572609
$data[14,5] = $null
$data[14,6] = "2025-10-01 08:54:34.766734+07:00"

$data[15,0] = "Telekomunikasi"
$data[15,1] = "https://www.cnnindonesia.com/teknologi/telekomunikasi"
$data[15,2] = "purbaya"
$data[15,3] = 1
$data[15,4] = 45931.37123720448
$data[15,5] = $null
$data[15,6] = "2025-10-01 08:54:34.894467+07:00"

$data[16,0] = "Climate"
$data[16,1] = "https://www.cnnindonesia.com/teknologi/climate"
$data[16,2] = "purbaya"
$data[16,3] = 1
$data[16,4] = 45931.3712387611
$data[16,5] = $null
$data[16,6] = "2025-10-01 08:54:35.028959+07:00"

$data[17,0] = "E-Vehicle"
$data[17,1] = "https://www.cnnindonesia.com/otomotif/e-vehicle"
$data[17,2] = "purbaya"
$data[17,3] = 1
$data[17,4] = 45931.37124012802
$data[17,5] = $null
$data[17,6] = "2025-10-01 08:54:35.147061+07:00"

$data[18,0] = "Commercial"
$data[18,1] = "https://www.cnnindonesia.com/otomotif/commercial"
$data[18,2] = "purbaya"
$data[18,3] = 1
$data[18,4] = 45931.37124148416
$data[18,5] = $null
$data[18,6] = "2025-10-01 08:54:35.264231+07:00"

$data[19,0] = "Info Otomotif"
$data[19,1] = "https://www.cnnindonesia.com/otomotif/info-otomotif"
$data[19,2] = "purbaya"
$data[19,3] = 1
$data[19,4] = 45931.37124283989
$data[19,5] = $null
$data[19,6] = "2025-10-01 08:54:35.381367+07:00"

$data[20,0] = "Seni Budaya"
$data[20,1] = "https://www.cnnindonesia.com/hiburan/seni-budaya"
$data[20,2] = "purbaya"
$data[20,3] = 1
$data[20,4] = 45931.37124359575
$data[20,5] = $null
$data[20,6] = "2025-10-01 08:54:35.446673+07:00"

$data[21,0] = "Health"
$data[21,1] = "https://www.cnnindonesia.com/gaya-hidup/health"
$data[21,2] = "purbaya"
$data[21,3] = 1
$data[21,4] = 45931.37124460488
$data[21,5] = $null
$data[21,6] = "2025-10-01 08:54:35.533862+07:00"

$data[22,0] = "Travel"
$data[22,1] = "https://www.cnnindonesia.com/gaya-hidup/travel"
$data[22,2] = "purbaya"
$data[22,3] = 1
$data[22,4] = 45931.37124763484
$data[22,5] = $null
$data[22,6] = "2025-10-01 08:54:35.795650+07:00"

$data[23,0] = "Trends"
$data[23,1] = "https://www.cnnindonesia.com/gaya-hidup/trends"
$data[23,2] = "purbaya"
$data[23,3] = 1
$data[23,4] = 45931.37124883582
$data[23,5] = $null
$data[23,6] = "2025-10-01 08:54:35.899415+07:00"

$data[24,0] = "Purbaya Respons Permintaan Bos BGN Tambah Anggaran MBG Rp28 T Ekonomi • 12 jam yang lalu"
$data[24,1] = "https://www.cnnindonesia.com/ekonomi/20250929180231-532-1278984/purbaya-respons-permintaan-bos-bgn-tambah-anggaran-mbg-rp28-t"
$data[24,2] = "purbaya"
$data[24,3] = 1
$data[24,4] = 45930.83377314815
$data[24,5] = "2025-09-30 20:00:38+07:00"
$data[24,6] = "2025-10-01 08:54:35.986398+07:00"

$data[25,0] = "Purbaya Sebut Ray Dalio Siap Dukung Program Indonesia Ekonomi • 13 jam yang lalu"
$data[25,1] = "https://www.cnnindonesia.com/ekonomi/20250930165101-532-1279378/purbaya-sebut-ray-dalio-siap-dukung-program-indonesia"
$data[25,2] = "purbaya"
$data[25,3] = 1
$data[25,4] = 45930.79541666667
$data[25,5] = "2025-09-30 19:05:24+07:00"
$data[25,6] = "2025-10-01 08:54:36.086780+07:00"

$data[26,0] = "Purbaya Sentil Balik Danantara yang Adukan Dirinya ke DPR Ekonomi • 15 jam yang lalu"
$data[26,1] = "https://www.cnnindonesia.com/ekonomi/20250930155921-532-1279361/purbaya-sentil-balik-danantara-yang-adukan-dirinya-ke-dpr"
$data[26,2] = "purbaya"
$data[26,3] = 1
$data[26,4] = 45930.73671296296
$data[26,5] = "2025-09-30 17:40:52+07:00"
$data[26,6] = "2025-10-01 08:54:36.157635+07:00"

$data[27,0] = "Purbaya Ditagih Subsidi BBM & Listrik: Sudah Bayar, Nyangkut di Mana? Ekonomi • 1 jam yang lalu"
$data[27,1] = "https://www.cnnindonesia.com/ekonomi/20250930113947-532-1279185/purbaya-ditagih-subsidi-bbm-listrik-sudah-bayar-nyangkut-di-mana"
$data[27,2] = "purbaya"
$data[27,3] = 1
$data[27,4] = 45930.51893518519
$data[27,5] = "2025-09-30 12:27:16+07:00"
$data[27,6] = "2025-10-01 08:54:36.254505+07:00"

$data[28,0] = "Purbaya Tiba-tiba Sidak Kantor BNI, Nimbrung Rapat Direksi Ekonomi • 23 jam yang lalu"
$data[28,1] = "https://www.cnnindonesia.com/ekonomi/20250929134914-532-1278863/purbaya-tiba-tiba-sidak-kantor-bni-nimbrung-rapat-direksi"
$data[28,2] = "purbaya"
$data[28,3] = 1
$data[28,4] = 45929.594375
$data[28,5] = "2025-09-29 14:15:54+07:00"
$data[28,6] = "2025-10-01 08:54:36.330606+07:00"

$data[29,0] = "Purbaya Ancam Ambil Alih Uang Pemda yang Masih Nganggur Ekonomi • 4 hari yang lalu"
$data[29,1] = "https://www.cnnindonesia.com/ekonomi/20250925180050-532-1277772/purbaya-ancam-ambil-alih-uang-pemda-yang-masih-nganggur"
$data[29,2] = "purbaya"
$data[29,3] = 1
$data[29,4] = 45925.81774305556
$data[29,5] = "2025-09-25 19:37:33+07:00"
$data[29,6] = "2025-10-01 08:54:36.406656+07:00"

$data[30,0] = "Luhut: Program MBG Bagus, Tinggal Pengelolaan yang Tertib Ekonomi • 4 hari yang lalu"
$data[30,1] = "https://www.cnnindonesia.com/ekonomi/20250925100133-532-1277517/luhut-program-mbg-bagus-tinggal-pengelolaan-yang-tertib"
$data[30,2] = "purbaya"
$data[30,3] = 1
$data[30,4] = 45925.744259259256
$data[30,5] = "2025-09-25 17:51:44+07:00"
$data[30,6] = "2025-10-01 08:54:39.418451+07:00"

$data[31,0] = "Purbaya Optimistis Mimpi Prabowo Defisit APBN 0 Persen Bisa Terwujud Ekonomi • 6 hari yang lalu"
$data[31,1] = "https://www.cnnindonesia.com/ekonomi/20250923163703-532-1276872/purbaya-optimistis-mimpi-prabowo-defisit-apbn-0-persen-bisa-terwujud"
$data[31,2] = "purbaya"
$data[31,3] = 1
$data[31,4] = 45923.71894675926
$data[31,5] = "2025-09-23 17:15:17+07:00"
$data[31,6] = "2025-10-01 08:54:41.837598+07:00"

$data[32,0] = "Purbaya Pede Bisa Genjot Ekonomi RI Lebih Cepat Tanpa Tambah Utang Ekonomi • 6 hari yang lalu"
$data[32,1] = "https://www.cnnindonesia.com/ekonomi/20250923142814-532-1276800/purbaya-pede-bisa-genjot-ekonomi-ri-lebih-cepat-tanpa-tambah-utang"
$data[32,2] = "purbaya"
$data[32,3] = 1
$data[32,4] = 45923.6375
$data[32,5] = "2025-09-23 15:18:00+07:00"
$data[32,6] = "2025-10-01 08:54:43.889120+07:00"

$data[33,0] = "Purbaya Tegaskan Data Ekonomi RI Tumbuh 5,12 Persen Bukan Manipulasi Ekonomi • 1 minggu yang lalu"
$data[33,1] = "https://www.cnnindonesia.com/ekonomi/20250922162321-532-1276432/purbaya-tegaskan-data-ekonomi-ri-tumbuh-512-persen-bukan-manipulasi"
$data[33,2] = "purbaya"
$data[33,3] = 1
$data[33,4] = 45922.812939814816
$data[33,5] = "2025-09-22 19:30:38+07:00"
$data[33,6] = "2025-10-01 08:54:44.025902+07:00"

$data[34,0] = "Purbaya Panggil Tokopedia Cs, Larang Jual Rokok Ilegal Mulai 1 Oktober Ekonomi • 1 minggu yang lalu"
$data[34,1] = "https://www.cnnindonesia.com/ekonomi/20250922174710-532-1276471/purbaya-panggil-tokopedia-cs-larang-jual-rokok-ilegal-mulai-1-oktober"
$data[34,2] = "purbaya"
$data[34,3] = 1
$data[34,4] = 45922.760983796295
$data[34,5] = "2025-09-22 18:15:49+07:00"
$data[34,6] = "2025-10-01 08:54:47.952041+07:00"

$data[35,0] = "Purbaya Ancam Sikat Oknum Bea Cukai Terlibat Peredaran Rokok Ilegal Ekonomi • 1 minggu yang lalu"
$data[35,1] = "https://www.cnnindonesia.com/ekonomi/20250922153320-532-1276410/purbaya-ancam-sikat-oknum-bea-cukai-terlibat-peredaran-rokok-ilegal"
$data[35,2] = "purbaya"
$data[35,3] = 1
$data[35,4] = 45922.68662037037
$data[35,5] = "2025-09-22 16:28:44+07:00"
$data[35,6] = "2025-10-01 08:54:48.074617+07:00"

$data[36,0] = "Menkeu Purbaya soal Satgas BLBI: Kalau Cuma Buat Keributan Enggak Usah Ekonomi • 1 minggu yang lalu"
$data[36,1] = "https://www.cnnindonesia.com/ekonomi/20250920072756-532-1275700/menkeu-purbaya-soal-satgas-blbi-kalau-cuma-buat-keributan-enggak-usah"
$data[36,2] = "purbaya"
$data[36,3] = 1
$data[36,4] = 45920.32638888889
$data[36,5] = "2025-09-20 07:50:00+07:00"
$data[36,6] = "2025-10-01 08:54:48.295212+07:00"

$data[37,0] = "Purbaya Tegas Tolak Tax Amnesty Jilid III: Beri Sinyal Boleh Melanggar Ekonomi • 1 minggu yang lalu"
$data[37,1] = "https://www.cnnindonesia.com/ekonomi/20250919140705-532-1275476/purbaya-tegas-tolak-tax-amnesty-jilid-iii-beri-sinyal-boleh-melanggar"
$data[37,2] = "purbaya"
$data[37,3] = 2
$data[37,4] = 45919.598287037035
$data[37,5] = "2025-09-19 14:21:32+07:00"
$data[37,6] = "2025-10-01 08:54:48.492464+07:00"

$data[38,0] = "Purbaya soal Gugatan Tutut Soeharto: Sudah Dicabut, Beliau Kirim Salam Ekonomi • 1 minggu yang lalu"
$data[38,1] = "https://www.cnnindonesia.com/ekonomi/20250918151932-532-1275123/purbaya-soal-gugatan-tutut-soeharto-sudah-dicabut-beliau-kirim-salam"
$data[38,2] = "purbaya"
$data[38,3] = 2
$data[38,4] = 45918.653657407405
$data[38,5] = "2025-09-18 15:41:16+07:00"
$data[38,6] = "2025-10-01 08:54:48.608264+07:00"

$data[39,0] = "Purbaya: Yang Takut RI Gelap, Jangan Khawatir, Sebentar Lagi Terang Ekonomi • 1 minggu yang lalu"
$data[39,1] = "https://www.cnnindonesia.com/ekonomi/20250917122423-532-1274591/purbaya-yang-takut-ri-gelap-jangan-khawatir-sebentar-lagi-terang"
$data[39,2] = "purbaya"
$data[39,3] = 2
$data[39,4] = 45917.52737268519
$data[39,5] = "2025-09-17 12:39:25+07:00"
$data[39,6] = "2025-10-01 08:54:48.737784+07:00"

$data[40,0] = "Purbaya Bersuara Saat Kebijakannya Dibandingkan dengan Sri Mulyani Ekonomi • 1 minggu yang lalu"
$data[40,1] = "https://www.cnnindonesia.com/ekonomi/20250917103601-532-1274537/purbaya-bersuara-saat-kebijakannya-dibandingkan-dengan-sri-mulyani"
$data[40,2] = "purbaya"
$data[40,3] = 2
$data[40,4] = 45917.451469907406
$data[40,5] = "2025-09-17 10:50:07+07:00"
$data[40,6] = "2025-10-01 08:54:48.860762+07:00"

$data[41,0] = "01:59 VIDEO: Purbaya: Dirut Bank Himbara Pusing, Saya Suruh Mereka Mikir TV • 2 minggu yang lalu"
$data[41,1] = "https://www.cnnindonesia.com/tv/20250916111325-400-1274175/video-purbaya-dirut-bank-himbara-pusing-saya-suruh-mereka-mikir"
$data[41,2] = "purbaya"
$data[41,3] = 2
$data[41,4] = 45916.46927083333
$data[41,5] = "2025-09-16 11:15:45+07:00"
$data[41,6] = "2025-10-01 08:54:48.973259+07:00"

$data[42,0] = "Purbaya Pede Defisit APBN Tak Melebar Gara-gara Stimulus Ekonomi Ekonomi • 2 minggu yang lalu"
$data[42,1] = "https://www.cnnindonesia.com/ekonomi/20250915175909-532-1273968/purbaya-pede-defisit-apbn-tak-melebar-gara-gara-stimulus-ekonomi"
$data[42,2] = "purbaya"
$data[42,3] = 2
$data[42,4] = 45915.87552083333
$data[42,5] = "2025-09-15 21:00:45+07:00"
$data[42,6] = "2025-10-01 08:54:49.091970+07:00"

$data[43,0] = "Purbaya Sebut Ada Bank Ngaku Tak Sanggup Serap Uang Pemerintah Rp200 T Ekonomi • 2 minggu yang lalu"
$data[43,1] = "https://www.cnnindonesia.com/ekonomi/20250915183126-532-1273985/purbaya-sebut-ada-bank-ngaku-tak-sanggup-serap-uang-pemerintah-rp200-t"
$data[43,2] = "purbaya"
$data[43,3] = 2
$data[43,4] = 45915.79498842593
$data[43,5] = "2025-09-15 19:04:47+07:00"
$data[43,6] = "2025-10-01 08:54:49.213245+07:00"

$data[44,0] = "Purbaya Guyur Bank Rp200 T, Zulhas Singgung Jatah Kopdes Merah Putih Ekonomi • 2 minggu yang lalu"
$data[44,1] = "https://www.cnnindonesia.com/ekonomi/20250915134736-92-1273846/purbaya-guyur-bank-rp200-t-zulhas-singgung-jatah-kopdes-merah-putih"
$data[44,2] = "purbaya"
$data[44,3] = 2
$data[44,4] = 45915.79400462963
$data[44,5] = "2025-09-15 19:03:22+07:00"
$data[44,6] = "2025-10-01 08:54:55.171422+07:00"

$data[45,0] = "05:08 VIDEO: Menanti Penyerapan Dana Suntikan Rp.200 Triliun TV • 2 minggu yang lalu"
$data[45,1] = "https://www.cnnindonesia.com/tv/20250914132720-404-1273533/video-menanti-penyerapan-dana-suntikan-rp200-triliun"
$data[45,2] = "purbaya"
$data[45,3] = 2
$data[45,4] = 45914.578668981485
$data[45,5] = "2025-09-14 13:53:17+07:00"
$data[45,6] = "2025-10-01 08:54:55.277981+07:00"

$data[46,0] = "Bank Himbara Dapat Suntikan Rp200 T, untuk Apa? Ekonomi • 2 minggu yang lalu"
$data[46,1] = "https://www.cnnindonesia.com/ekonomi/20250912173110-532-1273102/bank-himbara-dapat-suntikan-rp200-t-untuk-apa"
$data[46,2] = "purbaya"
$data[46,3] = 2
$data[46,4] = 45913.34037037037
$data[46,5] = "2025-09-13 08:10:08+07:00"
$data[46,6] = "2025-10-01 08:54:55.376215+07:00"

$data[47,0] = "Purbaya Akan Naikkan Jatah Anggaran Daerah Usai PBB Naik Selangit Ekonomi • 2 minggu yang lalu"
$data[47,1] = "https://www.cnnindonesia.com/ekonomi/20250911185023-532-1272721/purbaya-akan-naikkan-jatah-anggaran-daerah-usai-pbb-naik-selangit"
$data[47,2] = "purbaya"
$data[47,3] = 3
$data[47,4] = 45911.807488425926
$data[47,5] = "2025-09-11 19:22:47+07:00"
$data[47,6] = "2025-10-01 08:54:55.499675+07:00"

$data[48,0] = "Apa Efek Jika Purbaya Sebar Dana Pemerintah Rp200 T di BI ke Bank? Ekonomi • 2 minggu yang lalu"
$data[48,1] = "https://www.cnnindonesia.com/ekonomi/20250911111354-532-1272488/apa-efek-jika-purbaya-sebar-dana-pemerintah-rp200-t-di-bi-ke-bank"
$data[48,2] = "purbaya"
$data[48,3] = 3
$data[48,4] = 45911.50540509259
$data[48,5] = "2025-09-11 12:07:47+07:00"
$data[48,6] = "2025-10-01 08:54:55.618439+07:00"

$data[49,0] = "04:19 VIDEO: Rapat Perdana Menkeu Purbaya Dicecar Komisi XI DPR TV • 2 minggu yang lalu"
$data[49,1] = "https://www.cnnindonesia.com/tv/20250911092357-400-1272431/video-rapat-perdana-menkeu-purbaya-dicecar-komisi-xi-dpr"
$data[49,2] = "purbaya"
$data[49,3] = 3
$data[49,4] = 45911.39340277778
$data[49,5] = "2025-09-11 09:26:30+07:00"
$data[49,6] = "2025-10-01 08:54:55.746134+07:00"

$data[50,0] = "Serapan Anggaran Rendah, Purbaya Minta Bos BGN Jumpa Pers Tiap Bulan Ekonomi • 2 minggu yang lalu"
$data[50,1] = "https://www.cnnindonesia.com/ekonomi/20250910180046-532-1272257/serapan-anggaran-rendah-purbaya-minta-bos-bgn-jumpa-pers-tiap-bulan"
$data[50,2] = "purbaya"
$data[50,3] = 3
$data[50,4] = 45910.84768518519
$data[50,5] = "2025-09-10 20:20:40+07:00"
$data[50,6] = "2025-10-01 08:54:55.884359+07:00"

$data[51,0] = "00:55 VIDEO: Purbaya Tanya Balik DPR soal Warisan Masalah Fiskal dan Moneter Ekonomi • 2 minggu yang lalu"
$data[51,1] = "https://www.cnnindonesia.com/ekonomi/20250910175630-536-1272255/video-purbaya-tanya-balik-dpr-soal-warisan-masalah-fiskal-dan-moneter"
$data[51,2] = "purbaya"
$data[51,3] = 3
$data[51,4] = 45910.76162037037
$data[51,5] = "2025-09-10 18:16:44+07:00"
$data[51,6] = "2025-10-01 08:54:56.004451+07:00"

$data[52,0] = "Purbaya Sebut Rp425 T Mengendap di BI: Makanya Orang Susah Cari Kerja Ekonomi • 2 minggu yang lalu"
$data[52,1] = "https://www.cnnindonesia.com/ekonomi/20250910160549-532-1272212/purbaya-sebut-rp425-t-mengendap-di-bi-makanya-orang-susah-cari-kerja"
$data[52,2] = "purbaya"
$data[52,3] = 3
$data[52,4] = 45910.70972222222
$data[52,5] = "2025-09-10 17:02:00+07:00"
$data[52,6] = "2025-10-01 08:54:56.127726+07:00"

$data[53,0] = "Keluar"
$data[53,1] = "https://connect.detik.com/oauth/signout?redirectUrl=https%3A%2F%2Fwww.cnnindonesia.com%2Ftag%2Fsidak-bni"
$data[53,2] = "sidak BNI"
$data[53,3] = 1
$data[53,4] = 45931.37149996044
$data[53,5] = $null
$data[53,6] = "2025-10-01 08:54:57.596582+07:00"

$data[54,0] = "DAFTAR"
$data[54,1] = "https://connect.detik.com/accounts/register?clientId=10027&redirectUrl=https%3A%2F%2Fwww.cnnindonesia.com%2Fauthorize&backURL=https%3A%2F%2Fwww.cnnindonesia.com%2Ftag%2Fsidak-bni&ui=apps&osType=ANDROID"
$data[54,2] = "sidak BNI"
$data[54,3] = 1
$data[54,4] = 45931.37150020822
$data[54,5] = $null
$data[54,6] = "2025-10-01 08:54:57.617990+07:00"

$data[55,0] = "Keluar"
$data[55,1] = "https://connect.detik.com/oauth/signout?redirectUrl=https%3A%2F%2Fwww.cnnindonesia.com%2Ftag%2Fsidak-bni%3Fpage%3D2"
$data[55,2] = "sidak BNI"
$data[55,3] = 2
$data[55,4] = 45931.37158180954
$data[55,5] = $null
$data[55,6] = "2025-10-01 08:55:04.668344+07:00"

$data[56,0] = "DAFTAR"
$data[56,1] = "https://connect.detik.com/accounts/register?clientId=10027&redirectUrl=https%3A%2F%2Fwww.cnnindonesia.com%2Fauthorize&backURL=https%3A%2F%2Fwww.cnnindonesia.com%2Ftag%2Fsidak-bni%3Fpage%3D2&ui=apps&osType=ANDROID"
$data[56,2] = "sidak BNI"
$data[56,3] = 2
$data[56,4] = 45931.37158211042
$data[56,5] = $null
$data[56,6] = "2025-10-01 08:55:04.694340+07:00"

$data[57,0] = "Bea Cukai Beber Modus Jual Rokok Ilegal: Dipalsukan Jadi Pakaian Dalam Ekonomi • 1 hari yang lalu"
$data[57,1] = "https://www.cnnindonesia.com/ekonomi/20250929131723-532-1278848/bea-cukai-beber-modus-jual-rokok-ilegal-dipalsukan-jadi-pakaian-dalam"
$data[57,2] = "rokok ilegal"
$data[57,3] = 1
$data[57,4] = 45929.81300925926
$data[57,5] = "2025-09-29 19:30:44+07:00"
$data[57,6] = "2025-10-01 08:55:13.914726+07:00"

$data[58,0] = "Luhut Dukung Gaya Koboi Purbaya Pelototi MBG dan Berangus Rokok Ilegal Ekonomi • 5 hari yang lalu"
$data[58,1] = "https://www.cnnindonesia.com/ekonomi/20250924205413-532-1277414/luhut-dukung-gaya-koboi-purbaya-pelototi-mbg-dan-berangus-rokok-ilegal"
$data[58,2] = "rokok ilegal"
$data[58,3] = 1
$data[58,4] = 45925.39604166667
$data[58,5] = "2025-09-25 09:30:18+07:00"
$data[58,6] = "2025-10-01 08:55:14.028605+07:00"

$data[59,0] = "01:36 VIDEO: Menkeu Purbaya Ancam Tindak Penjual Rokok Ilegal di Marketplace TV • 6 hari yang lalu"
$data[59,1] = "https://www.cnnindonesia.com/tv/20250924133357-402-1277216/video-menkeu-purbaya-ancam-tindak-penjual-rokok-ilegal-di-marketplace"
$data[59,2] = "rokok ilegal"
$data[59,3] = 1
$data[59,4] = 45924.59372685185
$data[59,5] = "2025-09-24 14:14:58+07:00"
$data[59,6] = "2025-10-01 08:55:14.116666+07:00"

$data[60,0] = "Purbaya Komentari Tarif Cukai Rokok: Firaun Lu, Banyak Banget Ekonomi • 1 minggu yang lalu"
$data[60,1] = "https://www.cnnindonesia.com/ekonomi/20250919171305-92-1275594/purbaya-komentari-tarif-cukai-rokok-firaun-lu-banyak-banget"
$data[60,2] = "rokok ilegal"
$data[60,3] = 1
$data[60,4] = 45919.83064814815
$data[60,5] = "2025-09-19 19:56:08+07:00"
$data[60,6] = "2025-10-01 08:55:14.431082+07:00"

$data[61,0] = "TAIPAN Susilo Wonowidjojo, Raja Sigaret Asal Kediri Berharta Rp179 T Ekonomi • 2 minggu yang lalu"
$data[61,1] = "https://www.cnnindonesia.com/ekonomi/20250914044244-92-1273444/susilo-wonowidjojo-raja-sigaret-asal-kediri-berharta-rp179-t"
$data[61,2] = "rokok ilegal"
$data[61,3] = 1
$data[61,4] = 45914.67914351852
$data[61,5] = "2025-09-14 16:17:58+07:00"
$data[61,6] = "2025-10-01 08:55:14.553678+07:00"

$data[62,0] = "Wagub Jatim Janji Cegah Peredaran Rokok Ilegal Ekonomi • 2 minggu yang lalu"
$data[62,1] = "https://www.cnnindonesia.com/ekonomi/20250913201405-92-1273412/wagub-jatim-janji-cegah-peredaran-rokok-ilegal"
$data[62,2] = "rokok ilegal"
$data[62,3] = 1
$data[62,4] = 45914.0709375
$data[62,5] = "2025-09-14 01:42:09+07:00"
$data[62,6] = "2025-10-01 08:55:14.750864+07:00"

$data[63,0] = "Bos Bea Cukai Klaim Cegah Bocor Rp3,9 T dari Rokok hingga HP Ilegal Ekonomi • 2 bulan yang lalu"
$data[63,1] = "https://www.cnnindonesia.com/ekonomi/20250714174852-532-1250597/bos-bea-cukai-klaim-cegah-bocor-rp39-t-dari-rokok-hingga-hp-ilegal"
$data[63,2] = "rokok ilegal"
$data[63,3] = 1
$data[63,4] = 45854.60618055556
$data[63,5] = "2025-07-16 14:32:54+07:00"
$data[63,6] = "2025-10-01 08:55:14.894535+07:00"

$data[64,0] = "01:20 VIDEO: TNI AL Gagalkan Penyelundupan Jutaan Bungkus Rokok Ilegal TV • 2 bulan yang lalu"
$data[64,1] = "https://www.cnnindonesia.com/tv/20250701133042-407-1245599/video-tni-al-gagalkan-penyelundupan-jutaan-bungkus-rokok-ilegal"
$data[64,2] = "rokok ilegal"
$data[64,3] = 1
$data[64,4] = 45839.64591435185
$data[64,5] = "2025-07-01 15:30:07+07:00"
$data[64,6] = "2025-10-01 08:55:15.003247+07:00"

$data[65,0] = "Industri Tembakau Lesu, Bea Cukai Bentuk Satgas Cegah Rokok Ilegal Ekonomi • 3 bulan yang lalu"
$data[65,1] = "https://www.cnnindonesia.com/ekonomi/20250617174519-532-1240773/industri-tembakau-lesu-bea-cukai-bentuk-satgas-cegah-rokok-ilegal"
$data[65,2] = "rokok ilegal"
$data[65,3] = 2
$data[65,4] = 45825.76248842593
$data[65,5] = "2025-06-17 18:17:59+07:00"
$data[65,6] = "2025-10-01 08:55:18.862942+07:00"

$data[66,0] = "01:15 VIDEO: Bea Cukai Amankan Ratusan Ribu Bungkus Rokok Ilegal TV • 3 bulan yang lalu"
$data[66,1] = "https://www.cnnindonesia.com/tv/20250609152726-407-1237859/video-bea-cukai-amankan-ratusan-ribu-bungkus-rokok-ilegal"
$data[66,2] = "rokok ilegal"
$data[66,3] = 2
$data[66,4] = 45817.6669675926
$data[66,5] = "2025-06-09 16:00:26+07:00"
$data[66,6] = "2025-10-01 08:55:18.962910+07:00"

$data[67,0] = "Bea Cukai dan TNI AL Gagalkan Pengiriman Rokok Ilegal Nasional • 4 bulan yang lalu"
$data[67,1] = "https://www.cnnindonesia.com/nasional/20250518142657-12-1230456/bea-cukai-dan-tni-al-gagalkan-pengiriman-rokok-ilegal"
$data[67,2] = "rokok ilegal"
$data[67,3] = 2
$data[67,4] = 45795.6471875
$data[67,5] = "2025-05-18 15:31:57+07:00"
$data[67,6] = "2025-10-01 08:55:19.273533+07:00"

$data[68,0] = "Peredaran Rokok Ilegal Masih Marak, 752 Juta Batang Diamankan 2024 Ekonomi • 5 bulan yang lalu"
$data[68,1] = "https://www.cnnindonesia.com/ekonomi/20250417075106-532-1219649/peredaran-rokok-ilegal-masih-marak-752-juta-batang-diamankan-2024"
$data[68,2] = "rokok ilegal"
$data[68,3] = 2
$data[68,4] = 45764.43641203704
$data[68,5] = "2025-04-17 10:28:26+07:00"
$data[68,6] = "2025-10-01 08:55:19.437268+07:00"

$data[69,0] = "Polisi Sebut Mobil BRV Lawan Arah di Tol Pekalongan Bawa Rokok Ilegal Nasional • 5 bulan yang lalu"
$data[69,1] = "https://www.cnnindonesia.com/nasional/20250414181731-12-1218695/polisi-sebut-mobil-brv-lawan-arah-di-tol-pekalongan-bawa-rokok-ilegal"
$data[69,2] = "rokok ilegal"
$data[69,3] = 2
$data[69,4] = 45761.76472222222
$data[69,5] = "2025-04-14 18:21:12+07:00"
$data[69,6] = "2025-10-01 08:55:19.638125+07:00"

$data[70,0] = "Bakamla Tangkap Kapal Nihil Awak Angkut 200 Bal Rokok Ilegal di Kepri Nasional • 7 bulan yang lalu"
$data[70,1] = "https://www.cnnindonesia.com/nasional/20250215183408-20-1198683/bakamla-tangkap-kapal-nihil-awak-angkut-200-bal-rokok-ilegal-di-kepri"
$data[70,2] = "rokok ilegal"
$data[70,3] = 2
$data[70,4] = 45703.78006944444
$data[70,5] = "2025-02-15 18:43:18+07:00"
$data[70,6] = "2025-10-01 08:55:20.014284+07:00"

$data[71,0] = "Bea Cukai Sita 438,94 Juta Rokok Ilegal per 4 Agustus Ekonomi • 1 tahun yang lalu"
$data[71,1] = "https://www.cnnindonesia.com/ekonomi/20240906130245-532-1141723/bea-cukai-sita-43894-juta-rokok-ilegal-per-4-agustus"
$data[71,2] = "rokok ilegal"
$data[71,3] = 2
$data[71,4] = 45541.82098379629
$data[71,5] = "2024-09-06 19:42:13+07:00"
$data[71,6] = "2025-10-01 08:55:20.139695+07:00"

$data[72,0] = "Bea Cukai Gagalkan Impor 16 Kontainer Rokok Ilegal dari UAE Ekonomi • 1 tahun yang lalu"
$data[72,1] = "https://www.cnnindonesia.com/ekonomi/20240807205257-532-1130501/bea-cukai-gagalkan-impor-16-kontainer-rokok-ilegal-dari-uae"
$data[72,2] = "rokok ilegal"
$data[72,3] = 2
$data[72,4] = 45512.65280092593
$data[72,5] = "2024-08-08 15:40:02+07:00"
$data[72,6] = "2025-10-01 08:55:20.281343+07:00"

$data[73,0] = "7 Persen Rokok di Indonesia Ilegal Ekonomi • 1 tahun yang lalu"
$data[73,1] = "https://www.cnnindonesia.com/ekonomi/20240529135335-92-1103343/7-persen-rokok-di-indonesia-ilegal"
$data[73,2] = "rokok ilegal"
$data[73,3] = 2
$data[73,4] = 45442.47527777778
$data[73,5] = "2024-05-30 11:24:24+07:00"
$data[73,6] = "2025-10-01 08:55:20.409771+07:00"

$data[74,0] = "15,8 Juta Batang Rokok Ilegal Disita Ditjen Bea Cukai Setiap Minggu Ekonomi • 2 tahun yang lalu"
$data[74,1] = "https://www.cnnindonesia.com/ekonomi/20230811140806-532-984892/158-juta-batang-rokok-ilegal-disita-ditjen-bea-cukai-setiap-minggu"
$data[74,2] = "rokok ilegal"
$data[74,3] = 2
$data[74,4] = 45149.73644675926
$data[74,5] = "2023-08-11 17:40:29+07:00"
$data[74,6] = "2025-10-01 08:55:20.625580+07:00"

$data[75,0] = "Rokok Ilegal Meningkat Tahun Ini, Kerugian Rp548 M Ekonomi • 2 tahun yang lalu"
$data[75,1] = "https://www.cnnindonesia.com/ekonomi/20221212143947-532-886353/rokok-ilegal-meningkat-tahun-ini-kerugian-rp548-m"
$data[75,2] = "rokok ilegal"
$data[75,3] = 3
$data[75,4] = 44907.67791666667
$data[75,5] = "2022-12-12 16:16:12+07:00"
$data[75,6] = "2025-10-01 08:55:24.286972+07:00"

$data[76,0] = "Di Balik Kepulan Asap, Lampung 'Digempur' Rokok Ilegal Ekonomi • 2 tahun yang lalu"
$data[76,1] = "https://www.cnnindonesia.com/ekonomi/20221130105259-92-880734/di-balik-kepulan-asap-lampung-digempur-rokok-ilegal"
$data[76,2] = "rokok ilegal"
$data[76,3] = 3
$data[76,4] = 44896.39430555556
$data[76,5] = "2022-12-01 09:27:48+07:00"
$data[76,6] = "2025-10-01 08:55:24.432670+07:00"

$data[77,0] = "Bea Cukai Sumbagbar Buka Suara Soal Peredaran Rokok Ilegal di Lampung Ekonomi • 2 tahun yang lalu"
$data[77,1] = "https://www.cnnindonesia.com/ekonomi/20221130134845-92-880850/bea-cukai-sumbagbar-buka-suara-soal-peredaran-rokok-ilegal-di-lampung"
$data[77,2] = "rokok ilegal"
$data[77,3] = 3
$data[77,4] = 44895.87429398148
$data[77,5] = "2022-11-30 20:58:59+07:00"
$data[77,6] = "2025-10-01 08:55:24.648292+07:00"

$data[78,0] = "Bea Cukai Musnahkan Rokok Ilegal dan Sex Toys Senilai Rp2,2 M Ekonomi • 2 tahun yang lalu"
$data[78,1] = "https://www.cnnindonesia.com/ekonomi/20221130131940-532-880830/bea-cukai-musnahkan-rokok-ilegal-dan-sex-toys-senilai-rp22-m"
$data[78,2] = "rokok ilegal"
$data[78,3] = 3
$data[78,4] = 44895.57192129629
$data[78,5] = "2022-11-30 13:43:34+07:00"
$data[78,6] = "2025-10-01 08:55:24.794121+07:00"

$data[79,0] = "Bea Cukai Tindak 18.659 Kasus Rokok Ilegal Rugikan Negara Rp407 M 2022 Ekonomi • 2 tahun yang lalu"
$data[79,1] = "https://www.cnnindonesia.com/ekonomi/20221104162545-532-869632/bea-cukai-tindak-18659-kasus-rokok-ilegal-rugikan-negara-rp407-m-2022"
$data[79,2] = "rokok ilegal"
$data[79,3] = 3
$data[79,4] = 44869.73346064815
$data[79,5] = "2022-11-04 17:36:11+07:00"
$data[79,6] = "2025-10-01 08:55:25.416846+07:00"

$data[80,0] = "Bea Cukai Sulsel Musnahkan 5,2 Juta Rokok Ilegal Ekonomi • 3 tahun yang lalu"
$data[80,1] = "https://www.cnnindonesia.com/ekonomi/20220712155248-532-820536/bea-cukai-sulsel-musnahkan-52-juta-rokok-ilegal"
$data[80,2] = "rokok ilegal"
$data[80,3] = 3
$data[80,4] = 44754.94849537037
$data[80,5] = "2022-07-12 22:45:50+07:00"
$data[80,6] = "2025-10-01 08:55:25.614512+07:00"

$data[81,0] = "Bea Cukai Sita 1 Juta Rokok Impor China Ilegal untuk Pekerja Morowali Ekonomi • 3 tahun yang lalu"
$data[81,1] = "https://www.cnnindonesia.com/ekonomi/20220323071202-92-774953/bea-cukai-sita-1-juta-rokok-impor-china-ilegal-untuk-pekerja-morowali"
$data[81,2] = "rokok ilegal"
$data[81,3] = 3
$data[81,4] = 44643.31049768518
$data[81,5] = "2022-03-23 07:27:07+07:00"
$data[81,6] = "2025-10-01 08:55:25.832533+07:00"

$data[82,0] = "Bea Cukai Musnahkan Barang Ilegal Senilai RP15,6 Miliar di Tangerang Ekonomi • 3 tahun yang lalu"
$data[82,1] = "https://www.cnnindonesia.com/ekonomi/20211222133024-532-737253/bea-cukai-musnahkan-barang-ilegal-senilai-rp156-miliar-di-tangerang"
$data[82,2] = "rokok ilegal"
$data[82,3] = 3
$data[82,4] = 44552.64771990741
$data[82,5] = "2021-12-22 15:32:43+07:00"
$data[82,6] = "2025-10-01 08:55:26.002054+07:00"

$data[83,0] = "Rokok Ilegal Hingga Sex Toys Rugikan Negara Rp2 Miliar Dimusnahkan Nasional • 3 tahun yang lalu"
$data[83,1] = "https://www.cnnindonesia.com/nasional/20211118114053-12-722860/rokok-ilegal-hingga-sex-toys-rugikan-negara-rp2-miliar-dimusnahkan"
$data[83,2] = "rokok ilegal"
$data[83,3] = 3
$data[83,4] = 44518.51324074074
$data[83,5] = "2021-11-18 12:19:04+07:00"
$data[83,6] = "2025-10-01 08:55:26.188165+07:00"

$data[84,0] = "Rokok Ilegal Masih Dominasi Penindakan Bea Cukai Ekonomi • 4 tahun yang lalu"
$data[84,1] = "https://www.cnnindonesia.com/ekonomi/20210923124242-532-698403/rokok-ilegal-masih-dominasi-penindakan-bea-cukai"
$data[84,2] = "rokok ilegal"
$data[84,3] = 3
$data[84,4] = 44462.57931712963
$data[84,5] = "2021-09-23 13:54:13+07:00"
$data[84,6] = "2025-10-01 08:55:26.341979+07:00"

$ws.Range("A2:G86").Value = $data

# Ensure new row 86 inherits the date/time number format used in column E
$ws.Range("E86").NumberFormat = $ws.Range("E85").NumberFormat

Write-Host "Applied edits to rows 2-86"